$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 1752
$ws.Range("G4").Value = 78
$ws.Range("F5").Value = 441
$ws.Range("F8").Value = 324
$ws.Range("F9").Value = 304
$ws.Range("F10").Value = 1689
$ws.Range("F11").Value = 338
$ws.Range("F12").Value = 1394
$ws.Range("F13").Value = 789
$ws.Range("F15").Value = 661
$ws.Range("F16").Value = 12642
$ws.Range("F17").Value = 12667
$ws.Range("F18").Value = 935
$ws.Range("F21").Value = 298
$ws.Range("F23").Value = 504
$ws.Range("F27").Value = 231
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 49
$ws.Range("F5").Value = 72
$ws.Range("F9").Value = 50
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 82
$ws.Range("F3").Value = 153
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 82
$ws.Range("F4").Value = 153
$ws.Range("F6").Value = 1752
$ws.Range("G6").Value = 78
$ws.Range("F7").Value = 441
$ws.Range("F12").Value = 324
$ws.Range("F13").Value = 49
$ws.Range("F14").Value = 304
$ws.Range("F15").Value = 1689
$ws.Range("F16").Value = 338
$ws.Range("F17").Value = 1394
$ws.Range("F18").Value = 789
$ws.Range("F20").Value = 72
$ws.Range("F21").Value = 661
$ws.Range("F22").Value = 12642
$ws.Range("F23").Value = 12667
$ws.Range("F24").Value = 935
$ws.Range("F27").Value = 298
$ws.Range("F29").Value = 504
$ws.Range("F36").Value = 50
$ws.Range("F37").Value = 231
